$d = $word.ActiveDocument

$replacements = @(
    @("90×27=", "41×27="),
    @("69×58=", "77×19="),
    @("91×29=", "40×87="),
    @("93×23=", "46×22="),
    @("39×83=", "96×90="),
    @("98×26=", "73×64="),
    @("12×70=", "89×61="),
    @("11×45=", "20×93="),
    @("76×30=", "96×69="),
    @("19×26=", "69×90="),
    @("43×54=", "24×46="),
    @("39×66=", "22×77="),
    @("64×69=", "80×25="),
    @("86×38=", "90×63="),
    @("31×27=", "89×77="),
    @("12×20=", "82×60="),
    @("50×90=", "76×31="),
    @("78×94=", "69×64="),
    @("47×62=", "92×39="),
    @("58×14=", "97×97="),
    @("61×95=", "43×77="),
    @("70×22=", "11×99="),
    @("20×35=", "36×15="),
    @("89×53=", "21×25="),
    @("99×89=", "39×41=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
